$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r2 = New-Object 'object[,]' 1,12
$r2[0,0] = 1.088543089319103
$r2[0,1] = 0
$r2[0,2] = 0.1024648592858632
$r2[0,3] = 0.1916222410166717
$r2[0,4] = 2.896413622951727
$r2[0,5] = 2.126286655168514
$r2[0,6] = 1.78044049873975
$r2[0,7] = 0
$r2[0,8] = 0.3213029890215902
$r2[0,9] = 0.4098702357845241
$r2[0,10] = 0.1469987310750653
$r2[0,11] = 0.2313032576792047
$ws.Range("B2:M2").Value = $r2

$r3 = New-Object 'object[,]' 1,12
$r3[0,0] = 1.073125921524081
$r3[0,1] = 0
$r3[0,2] = 0.1026034800748175
$r3[0,3] = 0.192072795755418
$r3[0,4] = 2.892956076515731
$r3[0,5] = 2.118284681597459
$r3[0,6] = 1.781827795790633
$r3[0,7] = 0
$r3[0,8] = 0.3220425281297581
$r3[0,9] = 0.3713369223548852
$r3[0,10] = 0.1424872537828747
$r3[0,11] = 0.226977454838142
$ws.Range("B3:M3").Value = $r3

$r4 = New-Object 'object[,]' 1,12
$r4[0,0] = 1.064172484730136
$r4[0,1] = 0
$r4[0,2] = 0.1027060244492368
$r4[0,3] = 0.1923689792755767
$r4[0,4] = 2.891936610791802
$r4[0,5] = 2.114281596143769
$r4[0,6] = 1.783292469462438
$r4[0,7] = 0
$r4[0,8] = 0.3225206203393374
$r4[0,9] = 0.3478457736306666
$r4[0,10] = 0.139796846245396
$r4[0,11] = 0.224436436706732
$ws.Range("B4:M4").Value = $r4

$r5 = New-Object 'object[,]' 1,12
$r5[0,0] = 1.060653250920211
$r5[0,1] = 0
$r5[0,2] = 0.1027522121569113
$r5[0,3] = 0.1924946002801686
$r5[0,4] = 2.891798760669474
$r5[0,5] = 2.112879094983754
$r5[0,6] = 1.784043536132785
$r5[0,7] = 0
$r5[0,8] = 0.3227214948317338
$r5[0,9] = 0.3383155567797473
$r5[0,10] = 0.1387205898297807
$r5[0,11] = 0.2234299792617271
$ws.Range("B5:M5").Value = $r5

$r6 = New-Object 'object[,]' 1,12
$r6[0,0] = 1.060076710677777
$r6[0,1] = 0
$r6[0,2] = 0.1027601478454763
$r6[0,3] = 0.192515757234839
$r6[0,4] = 2.891792639852014
$r6[0,5] = 2.112660026763464
$r6[0,6] = 1.784177566159784
$r6[0,7] = 0
$r6[0,8] = 0.3227552154679749
$r6[0,9] = 0.3367356531494465
$r6[0,10] = 0.1385430949077318
$r6[0,11] = 0.223264613945549
$ws.Range("B6:M6").Value = $r6

$r7 = New-Object 'object[,]' 1,12
$r7[0,0] = 1.064124498755461
$r7[0,1] = 0
$r7[0,2] = 0.102706629513829
$r7[0,3] = 0.1923706534942498
$r7[0,4] = 2.891933627604004
$r7[0,5] = 2.114261755275606
$r7[0,6] = 1.783301974146823
$r7[0,7] = 0
$r7[0,8] = 0.3225233049013521
$r7[0,9] = 0.347717072895108
$r7[0,10] = 0.1397822499628063
$r7[0,11] = 0.2244227456031638
$ws.Range("B7:M7").Value = $r7

$r8 = New-Object 'object[,]' 1,12
$r8[0,0] = 1.083121071950472
$r8[0,1] = 0
$r8[0,2] = 0.1025090480609343
$r8[0,3] = 0.1917735425784666
$r8[0,4] = 2.894992492841496
$r8[0,5] = 2.123338617177481
$r8[0,6] = 1.780791688847813
$r8[0,7] = 0
$r8[0,8] = 0.3215530044268959
$r8[0,9] = 0.3965491248623039
$r8[0,10] = 0.1454266842226914
$r8[0,11] = 0.2297879015144524
$ws.Range("B8:M8").Value = $r8

$r9 = New-Object 'object[,]' 1,12
$r9[0,0] = 1.124426390367091
$r9[0,1] = 0
$r9[0,2] = 0.1022591823264207
$r9[0,3] = 0.1907572104125612
$r9[0,4] = 2.909744841936131
$r9[0,5] = 2.148366739522061
$r9[0,6] = 1.780729273832378
$r9[0,7] = 0
$r9[0,8] = 0.3198403038209312
$r9[0,9] = 0.4936394253541039
$r9[0,10] = 0.1571248984757574
$r9[0,11] = 0.2412182163482797
$ws.Range("B9:M9").Value = $r9

$r10 = New-Object 'object[,]' 1,12
$r10[0,0] = 1.157228203055183
$r10[0,1] = 0
$r10[0,2] = 0.1021585726943606
$r10[0,3] = 0.1901041444087994
$r10[0,4] = 2.925922425286572
$r10[0,5] = 2.171175471600691
$r10[0,6] = 1.783644898945766
$r10[0,7] = 0
$r10[0,8] = 0.3186971653873583
$r10[0,9] = 0.5657830133987432
$r10[0,10] = 0.1661010586913818
$r10[0,11] = 0.2501668398650239
$ws.Range("B10:M10").Value = $r10

$r11 = New-Object 'object[,]' 1,12
$r11[0,0] = 1.172680562524505
$r11[0,1] = 0
$r11[0,2] = 0.1021306314585431
$r11[0,3] = 0.1898272497070552
$r11[0,4] = 2.934442200399047
$r11[0,5] = 2.182514938971167
$r11[0,6] = 1.785614201243874
$r11[0,7] = 0
$r11[0,8] = 0.3182019875749802
$r11[0,9] = 0.5987797470906742
$r11[0,10] = 0.170266985656454
$r11[0,11] = 0.2543567347253912
$ws.Range("B11:M11").Value = $r11

$r12 = New-Object 'object[,]' 1,12
$r12[0,0] = 1.178607925995351
$r12[0,1] = 0
$r12[0,2] = 0.102122598642346
$r12[0,3] = 0.1897252899877975
$r12[0,4] = 2.937835280194491
$r12[0,5] = 2.186947650202569
$r12[0,6] = 1.786452336894229
$r12[0,7] = 0
$r12[0,8] = 0.3180180385142037
$r12[0,9] = 0.6113002603852067
$r12[0,10] = 0.1718563392991825
$r12[0,11] = 0.2559603887376412
$ws.Range("B12:M12").Value = $r12

$r13 = New-Object 'object[,]' 1,12
$r13[0,0] = 1.177327993581599
$r13[0,1] = 0
$r13[0,2] = 0.1021242155778452
$r13[0,3] = 0.1897471202517629
$r13[0,4] = 2.937097100980338
$r13[0,5] = 2.185986815725244
$r13[0,6] = 1.786267720569583
$r13[0,7] = 0
$r13[0,8] = 0.318057496853295
$r13[0,9] = 0.6086026198479715
$r13[0,10] = 0.1715135195834279
$r13[0,11] = 0.2556142572086699
$ws.Range("B13:M13").Value = $r13

$r14 = New-Object 'object[,]' 1,12
$r14[0,0] = 1.173166691706342
$r14[0,1] = 0
$r14[0,2] = 0.1021299196179264
$r14[0,3] = 0.1898188034634885
$r14[0,4] = 2.934718007986078
$r14[0,5] = 2.182876840284962
$r14[0,6] = 1.785681303572915
$r14[0,7] = 0
$r14[0,8] = 0.3181867825925915
$r14[0,9] = 0.5998093103860072
$r14[0,10] = 0.170397506518384
$r14[0,11] = 0.2544883275212015
$ws.Range("B14:M14").Value = $r14

$r15 = New-Object 'object[,]' 1,12
$r15[0,0] = 1.17062764499147
$r15[0,1] = 0
$r15[0,2] = 0.1021337448651778
$r15[0,3] = 0.1898630881647057
$r15[0,4] = 2.933282468884059
$r15[0,5] = 2.180989957290222
$r15[0,6] = 1.785334138298367
$r15[0,7] = 0
$r15[0,8] = 0.3182664377421296
$r15[0,9] = 0.5944264500529073
$r15[0,10] = 0.1697154517947581
$r15[0,11] = 0.2538008782283399
$ws.Range("B15:M15").Value = $r15

$r16 = New-Object 'object[,]' 1,12
$r16[0,0] = 1.156228997395687
$r16[0,1] = 0
$r16[0,2] = 0.1021607561041442
$r16[0,3] = 0.1901226456185059
$r16[0,4] = 2.925388985748185
$r16[0,5] = 2.170453812982657
$r16[0,6] = 1.78352913391268
$r16[0,7] = 0
$r16[0,8] = 0.3187300255982146
$r16[0,9] = 0.5636301633163043
$r16[0,10] = 0.1658304618721473
$r16[0,11] = 0.2498954090336696
$ws.Range("B16:M16").Value = $r16

$r17 = New-Object 'object[,]' 1,12
$r17[0,0] = 1.147531508522349
$r17[0,1] = 0
$r17[0,2] = 0.1021818813224371
$r17[0,3] = 0.1902870401586627
$r17[0,4] = 2.920843785512574
$r17[0,5] = 2.164237137299295
$r17[0,6] = 1.782586463023677
$r17[0,7] = 0
$r17[0,8] = 0.3190207784247741
$r17[0,9] = 0.5447831032038266
$r17[0,10] = 0.1634682560594314
$r17[0,11] = 0.2475299701931633
$ws.Range("B17:M17").Value = $r17

$r18 = New-Object 'object[,]' 1,12
$r18[0,0] = 1.142578915021858
$r18[0,1] = 0
$r18[0,2] = 0.1021957108266172
$r18[0,3] = 0.1903834962881037
$r18[0,4] = 2.918338742547718
$r18[0,5] = 2.160752166526265
$r18[0,6] = 1.782104786213722
$r18[0,7] = 0
$r18[0,8] = 0.3191903504520206
$r18[0,9] = 0.5339595664618173
$r18[0,10] = 0.1621173593421474
$r18[0,11] = 0.2461806503738906
$ws.Range("B18:M18").Value = $r18

$r19 = New-Object 'object[,]' 1,12
$r19[0,0] = 1.140910647446105
$r19[0,1] = 0
$r19[0,2] = 0.1022006820775978
$r19[0,3] = 0.1904164813975617
$r19[0,4] = 2.917509341103454
$r19[0,5] = 2.159587788653909
$r19[0,6] = 1.781952095217179
$r19[0,7] = 0
$r19[0,8] = 0.3192481665490283
$r19[0,9] = 0.53029779541518
$r19[0,10] = 0.1616613076922704
$r19[0,11] = 0.2457257235665438
$ws.Range("B19:M19").Value = $r19

$r20 = New-Object 'object[,]' 1,12
$r20[0,0] = 1.148452202165856
$r20[0,1] = 0
$r20[0,2] = 0.1021794588675995
$r20[0,3] = 0.1902693434241338
$r20[0,4] = 2.921316324708897
$r20[0,5] = 2.16488952534371
$r20[0,6] = 1.782680548484194
$r20[0,7] = 0
$r20[0,8] = 0.3189895853040148
$r20[0,9] = 0.5467876687808371
$r20[0,10] = 0.1637189120759075
$r20[0,11] = 0.2477806148813499
$ws.Range("B20:M20").Value = $r20

$r21 = New-Object 'object[,]' 1,12
$r21[0,0] = 1.174386910603943
$r21[0,1] = 0
$r21[0,2] = 0.1021281751758423
$r21[0,3] = 0.1897976698855608
$r21[0,4] = 2.935412278091206
$r21[0,5] = 2.18378655013359
$r21[0,6] = 1.785851041083276
$r21[0,7] = 0
$r21[0,8] = 0.3181487115746298
$r21[0,9] = 0.6023914321947927
$r21[0,10] = 0.1707249867421439
$r21[0,11] = 0.2548185789011512
$ws.Range("B21:M21").Value = $r21

$r22 = New-Object 'object[,]' 1,12
$r22[0,0] = 1.19177889709178
$r22[0,1] = 0
$r22[0,2] = 0.1021095025341729
$r22[0,3] = 0.1895062701212513
$r22[0,4] = 2.945597162603917
$r22[0,5] = 2.196945363352029
$r22[0,6] = 1.78846171548696
$r22[0,7] = 0
$r22[0,8] = 0.3176199221571649
$r22[0,9] = 0.6388794397707898
$r22[0,10] = 0.1753726509991651
$r22[0,11] = 0.2595175245779657
$ws.Range("B22:M22").Value = $r22

$r23 = New-Object 'object[,]' 1,12
$r23[0,0] = 1.18245614630888
$r23[0,1] = 0
$r23[0,2] = 0.1021181152919439
$r23[0,3] = 0.1896602552765083
$r23[0,4] = 2.940072344439329
$r23[0,5] = 2.189848235680643
$r23[0,6] = 1.787019087695768
$r23[0,7] = 0
$r23[0,8] = 0.3179002492214553
$r23[0,9] = 0.6193916727218038
$r23[0,10] = 0.1728858351330729
$r23[0,11] = 0.257000561430722
$ws.Range("B23:M23").Value = $r23

$r24 = New-Object 'object[,]' 1,12
$r24[0,0] = 1.148035807930967
$r24[0,1] = 0
$r24[0,2] = 0.1021805488121643
$r24[0,3] = 0.1902773380639227
$r24[0,4] = 2.921102353143183
$r24[0,5] = 2.164594303235134
$r24[0,6] = 1.782637824671895
$r24[0,7] = 0
$r24[0,8] = 0.3190036801932274
$r24[0,9] = 0.5458813676330294
$r24[0,10] = 0.1636055681618274
$r24[0,11] = 0.2476672653849761
$ws.Range("B24:M24").Value = $r24

$r25 = New-Object 'object[,]' 1,12
$r25[0,0] = 1.112819968212051
$r25[0,1] = 0
$r25[0,2] = 0.1023121427905274
$r25[0,3] = 0.1910156680558797
$r25[0,4] = 2.904816527094226
$r25[0,5] = 2.140820677271932
$r25[0,6] = 1.780225989753859
$r25[0,7] = 0
$r25[0,8] = 0.3202833586094624
$r25[0,9] = 0.4672314656159244
$r25[0,10] = 0.1538930682806381
$r25[0,11] = 0.2380290395334228
$ws.Range("B25:M25").Value = $r25
